# Weekly update: insert the latest week's "Pera" (Packham's Triumph)
# price records at the top of the data block (row 501), pushing the
# existing history down by two rows (dimension grows from T602 to T604).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two blank rows at 501:502 - Excel shifts rows 501:602 down to
# 503:604 and carries the column D date-number formatting with them.
$ws.Rows("501:502").Insert()

# New row 501 - "Primera" grade
$ws.Cells.Item(501, 1).Value = 4
$ws.Cells.Item(501, 2).Value = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item(501, 3).Value = "Los Lagos"
$ws.Cells.Item(501, 4).Value = 45211
$ws.Cells.Item(501, 5).Value = 10
$ws.Cells.Item(501, 6).Value = "Fruta"
$ws.Cells.Item(501, 7).Value = 100104
$ws.Cells.Item(501, 8).Value = "Frutos de pepita"
$ws.Cells.Item(501, 9).Value = 100104005
$ws.Cells.Item(501, 10).Value = "Pera"
$ws.Cells.Item(501, 11).Value = "Packham's Triumph"
$ws.Cells.Item(501, 12).Value = "Primera"
$ws.Cells.Item(501, 13).Value = 200
$ws.Cells.Item(501, 14).Value = 20000
$ws.Cells.Item(501, 15).Value = 20000
$ws.Cells.Item(501, 16).Value = 20000
$ws.Cells.Item(501, 17).Value = "`$/caja 15 kilos empedrada"
$ws.Cells.Item(501, 18).Value = "Región de O'Higgins"
$ws.Cells.Item(501, 19).Value = 1333
$ws.Cells.Item(501, 20).Value = 15

# New row 502 - "Segunda" grade
$ws.Cells.Item(502, 1).Value = 4
$ws.Cells.Item(502, 2).Value = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item(502, 3).Value = "Los Lagos"
$ws.Cells.Item(502, 4).Value = 45211
$ws.Cells.Item(502, 5).Value = 10
$ws.Cells.Item(502, 6).Value = "Fruta"
$ws.Cells.Item(502, 7).Value = 100104
$ws.Cells.Item(502, 8).Value = "Frutos de pepita"
$ws.Cells.Item(502, 9).Value = 100104005
$ws.Cells.Item(502, 10).Value = "Pera"
$ws.Cells.Item(502, 11).Value = "Packham's Triumph"
$ws.Cells.Item(502, 12).Value = "Segunda"
$ws.Cells.Item(502, 13).Value = 200
$ws.Cells.Item(502, 14).Value = 16000
$ws.Cells.Item(502, 15).Value = 16000
$ws.Cells.Item(502, 16).Value = 16000
$ws.Cells.Item(502, 17).Value = "`$/caja 15 kilos empedrada"
$ws.Cells.Item(502, 18).Value = "Región de O'Higgins"
$ws.Cells.Item(502, 19).Value = 1067
$ws.Cells.Item(502, 20).Value = 15
